# Applies the "Added a few more slots" change to the Casino Island II
# review document:
#
#   1. Insert a new paragraph right after the H1 title, containing a bold
#      "Meta description" run followed by a plain run with the SEO blurb
#      (this is the text that used to live at the very end of the doc).
#
#   2. Near the end of the document, delete the paragraph that duplicated
#      the bold title text ("Play Casino Island II Slot Game for Free -
#      Review"), and rewrite the following italic paragraph's text with
#      the new feature-image prompt (keeping the italic run formatting
#      and the leading empty run untouched).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert the "Meta description" paragraph after the H1 title.
# ---------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)

# Build the new paragraph's contents (leading empty run + bold label run +
# plain run) directly as OOXML so the formatting/run split matches exactly,
# then drop it into the freshly-inserted (still empty) paragraph.
$metaXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Discover Casino Island II, the slot machine featuring a versatile range of bets, a bonus game, and an impressive RTP of 96.9%. Play for free now!</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$null = $metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# Part 2: near the end, drop the duplicated bold title paragraph and
# rewrite the italic paragraph's text.
# ---------------------------------------------------------------------

$oldTitleText = "Play Casino Island II Slot Game for Free - Review"
$oldMetaText = "Discover Casino Island II, the slot machine featuring a versatile range of bets, a bonus game, and an impressive RTP of 96.9%. Play for free now!"

$boldTitlePara = $null
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq $oldTitleText -and $para.Range.Font.Bold) {
        $boldTitlePara = $para
        break
    }
}

if ($boldTitlePara -ne $null) {
    $boldTitlePara.Range.Delete()
}

$italicPara = $null
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq $oldMetaText) {
        $italicPara = $para
        break
    }
}

$newImagePromptText = "Create a fun and adventurous feature image for `"Casino Island II`" that highlights the game's seafaring theme and unique gameplay. The image should be in cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should look like they're enjoying the game and have found treasure on the island in the background. The image should also include some of the game's symbols, such as the sailboat, helm, lighthouse, treasure chest, and compass, to emphasize the game's seafaring elements. Use bright and vibrant colors to make the image pop and convey a sense of excitement and adventure to potential players."

if ($italicPara -ne $null) {
    $fullRange = $italicPara.Range
    # Exclude the trailing paragraph mark so only the visible text run is
    # replaced - this keeps the leading empty run and the <w:i/> run
    # formatting intact.
    $textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
    $textRange.Text = $newImagePromptText
}
